$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to Text format before writing, so that
# values like "1.000", "2.370", "23.066.20" are stored as literal text
# (verbatim strings) instead of being auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '23.066.20'
$ws.Range("D3").Value = '1.602.70'
$ws.Range("E3").Value = '  -2.94%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '1.002'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = '301.24'
$ws.Range("E6").Value = '  -3.14%  '
$ws.Range("D7").Value = '0.3778'
$ws.Range("E7").Value = '  -3.02%  '
$ws.Range("D8").Value = '0.3653'
$ws.Range("E8").Value = '  -4.46%  '
$ws.Range("D9").Value = '50.08'
$ws.Range("E9").Value = '  -3.35%  '
$ws.Range("D10").Value = '1.266'
$ws.Range("E10").Value = '  -5.68%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '0.08139'
$ws.Range("E12").Value = '  -3.71%  '
$ws.Range("D13").Value = '22.83'
$ws.Range("D14").Value = '6.605'
$ws.Range("E14").Value = '  -5.84%  '
$ws.Range("D15").Value = '0.00001257'
$ws.Range("E15").Value = '  -4.35%  '
$ws.Range("D16").Value = '7.389'
$ws.Range("E16").Value = '  -7.95%  '
$ws.Range("D17").Value = '1.601.41'
$ws.Range("E17").Value = '  -2.96%  '
$ws.Range("D18").Value = '91.97'
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("D19").Value = '0.06881'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").Value = '18.24'
$ws.Range("E20").Value = '  -6.79%  '
$ws.Range("D21").Value = '6.588'
$ws.Range("E21").Value = '  -5.38%  '
$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").Value = '0.5573'
$ws.Range("E22").Value = '  -5.49%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '13.11'
$ws.Range("E24").Value = '  -4.11%  '
$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").Value = '23.096.34'
$ws.Range("E25").Value = '  -3.32%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '2.354'
$ws.Range("E26").Value = '  -3.84%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.781'
$ws.Range("E27").Value = '  -5.19%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '21.16'
$ws.Range("E28").Value = '  -3.99%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '150.36'
$ws.Range("E29").Value = '  -2.04%  '
$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").Value = '5.266'
$ws.Range("E30").Value = '  -2.73%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '133.47'
$ws.Range("E31").Value = '  -2.85%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '2.370'
$ws.Range("E32").Value = '  -4.63%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.817'
$ws.Range("E33").Value = '  -12.06%  '
$ws.Range("B34").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C34").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D34").Value = '1.775.75'
$ws.Range("E34").Value = '  -3.11%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.9530'
$ws.Range("E35").Value = '  -4.64%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.07678'
$ws.Range("E36").Value = '  -5.60%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '10.43'
$ws.Range("E37").Value = '  -2.61%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '6.285'
$ws.Range("E38").Value = '  -5.69%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.02724'
$ws.Range("E39").Value = '  -6.17%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.2545'
$ws.Range("E40").Value = '  -4.88%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.08909'
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.369'
$ws.Range("E42").Value = '  -3.65%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.7096'
$ws.Range("E43").Value = '  -6.05%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '12.65'
$ws.Range("E44").Value = '  -6.41%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '15.37'
$ws.Range("E45").Value = '  -6.42%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6627'
$ws.Range("E46").Value = '  -4.53%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '2.332'
$ws.Range("E47").Value = '  -4.44%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '4.003'
$ws.Range("E49").Value = '  -2.54%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '132.51'
$ws.Range("E50").Value = '  -0.80%  '
$ws.Range("D51").Value = '1.247'
$ws.Range("E51").Value = '  +1.55%  '

# Restore the default (unstyled) cell style on column D now that the
# literal text values are committed, matching the original workbook
# formatting (data rows carry no explicit style).
$ws.Range("D2:D51").Style = "Normal"
